# Auto commit at 2025-11-10 8:19:09.26
#
# Refresh the "Metrics" sheet's raw indicator values (B2:B13) with the
# latest pull, which ripples into the formula-driven "today" dashboard
# sheet (B11:B22 / E11:E22 / F11:F22 all reference Metrics!B2:B13, and
# A1 = TODAY()-1 is recalculated automatically). Also nudges the saved
# cursor position on both sheets to where the author last clicked.

$wb = $excel.ActiveWorkbook

# ---- Metrics sheet: update the raw metric values ----------------------
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 119303.91
$wsMetrics.Range("B3").Value  = 104268.17
$wsMetrics.Range("B4").Value  = 37079.58
$wsMetrics.Range("B5").Value  = 4968
$wsMetrics.Range("B6").Value  = 4915549.66
$wsMetrics.Range("B7").Value  = 4146344.850000001
$wsMetrics.Range("B8").Value  = 1444039.41
$wsMetrics.Range("B9").Value  = 191175
$wsMetrics.Range("B10").Value = 33380930.650000006
$wsMetrics.Range("B11").Value = 31421620.009999998
$wsMetrics.Range("B12").Value = 11725761.450000003
$wsMetrics.Range("B13").Value = 1288805

# Saved selection on Metrics moved from D13 to D7. Activate the sheet
# just long enough to move the selection, the same way the original
# editing session would have (we re-activate "today" below so it stays
# the sheet that's on top when the workbook is saved).
[void]$wsMetrics.Activate()
[void]$wsMetrics.Range("D7").Select()

# ---- today sheet: just move the saved selection ------------------------
# B11:B22 / E11:E22 / F11:F22 and A1 (TODAY()-1) are all formulas, so
# they recalc automatically from the Metrics change above - nothing to
# write here directly.
$wsToday = $wb.Worksheets.Item("today")
[void]$wsToday.Activate()
[void]$wsToday.Range("H12").Select()
